{"js": "// Apply the multiplication-table refresh + date bump described in the diff.\n// Each (old, new) pair is unique in the document, so body.search + insertText(..., replace)\n// on every match safely targets exactly the run that needs updating.\nconst replacements = [\n  [\"2025-08-05 Tuesday\", \"2025-08-06 Wednesday\"],\n  [\"81\u00d767=5427\", \"16\u00d788=1408\"],\n  [\"17\u00d764=1088\", \"24\u00d719=456\"],\n  [\"22\u00d711=242\", \"15\u00d729=435\"],\n  [\"83\u00d747=3901\", \"92\u00d712=1104\"],\n  [\"25\u00d763=1575\", \"97\u00d759=5723\"],\n  [\"65\u00d730=1950\", \"93\u00d773=6789\"],\n  [\"40\u00d792=3680\", \"99\u00d761=6039\"],\n  [\"42\u00d764=2688\", \"99\u00d769=6831\"],\n  [\"48\u00d795=4560\", \"20\u00d732=640\"],\n  [\"67\u00d761=4087\", \"16\u00d750=800\"],\n  [\"15\u00d732=480\", \"20\u00d788=1760\"],\n  [\"20\u00d730=600\", \"69\u00d742=2898\"],\n  [\"40\u00d746=1840\", \"65\u00d715=975\"],\n  [\"87\u00d771=6177\", \"26\u00d788=2288\"],\n  [\"16\u00d785=1360\", \"34\u00d718=612\"],\n  [\"91\u00d787=7917\", \"12\u00d717=204\"],\n  [\"32\u00d749=1568\", \"79\u00d725=1975\"],\n  [\"69\u00d730=2070\", \"67\u00d736=2412\"],\n  [\"42\u00d784=3528\", \"97\u00d799=9603\"],\n  [\"84\u00d720=1680\", \"47\u00d761=2867\"],\n  [\"51\u00d763=3213\", \"77\u00d711=847\"],\n  [\"29\u00d713=377\", \"39\u00d792=3588\"],\n  [\"75\u00d763=4725\", \"72\u00d749=3528\"],\n  [\"49\u00d785=4165\", \"72\u00d754=3888\"],\n  [\"34\u00d786=2924\", \"56\u00d733=1848\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Apply the multiplication-table refresh + date bump described in the diff.\n# Each (old, new) pair is unique in the document, so Find/Replace (wdReplaceAll)\n# on each pair safely retargets exactly the run that needs updating, leaving\n# all run/paragraph formatting (fonts, size, alignment) untouched.\n$replacements = @(\n    @(\"2025-08-05 Tuesday\", \"2025-08-06 Wednesday\"),\n    @(\"81\u00d767=5427\", \"16\u00d788=1408\"),\n    @(\"17\u00d764=1088\", \"24\u00d719=456\"),\n    @(\"22\u00d711=242\", \"15\u00d729=435\"),\n    @(\"83\u00d747=3901\", \"92\u00d712=1104\"),\n    @(\"25\u00d763=1575\", \"97\u00d759=5723\"),\n    @(\"65\u00d730=1950\", \"93\u00d773=6789\"),\n    @(\"40\u00d792=3680\", \"99\u00d761=6039\"),\n    @(\"42\u00d764=2688\", \"99\u00d769=6831\"),\n    @(\"48\u00d795=4560\", \"20\u00d732=640\"),\n    @(\"67\u00d761=4087\", \"16\u00d750=800\"),\n    @(\"15\u00d732=480\", \"20\u00d788=1760\"),\n    @(\"20\u00d730=600\", \"69\u00d742=2898\"),\n    @(\"40\u00d746=1840\", \"65\u00d715=975\"),\n    @(\"87\u00d771=6177\", \"26\u00d788=2288\"),\n    @(\"16\u00d785=1360\", \"34\u00d718=612\"),\n    @(\"91\u00d787=7917\", \"12\u00d717=204\"),\n    @(\"32\u00d749=1568\", \"79\u00d725=1975\"),\n    @(\"69\u00d730=2070\", \"67\u00d736=2412\"),\n    @(\"42\u00d784=3528\", \"97\u00d799=9603\"),\n    @(\"84\u00d720=1680\", \"47\u00d761=2867\"),\n    @(\"51\u00d763=3213\", \"77\u00d711=847\"),\n    @(\"29\u00d713=377\", \"39\u00d792=3588\"),\n    @(\"75\u00d763=4725\", \"72\u00d749=3528\"),\n    @(\"49\u00d785=4165\", \"72\u00d754=3888\"),\n    @(\"34\u00d786=2924\", \"56\u00d733=1848\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $result) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\n"}
